$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Word count: 4,900 -> ~5,000
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "4,900 words (main text, excluding abstract, references, tables, and figure legends)",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "~5,000 words (main text, excluding abstract, references, tables, and figure legends)", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2. "Number of Figures:" -> "Number of Main Figures:" (label text only,
#    keeps the existing bold run/formatting intact since the whole match
#    sits inside a single run)
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Number of Figures:", $true, $false, $false, $false, $false, $true, 1, $false,
    "Number of Main Figures:", 2) | Out-Null

# "3 main figures" -> "3"
$d.Content.Find.Execute("3 main figures", $true, $false, $false, $false, $false, $true, 1, $false,
    "3", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3. Insert the brand-new "Number of Supplementary Figures: 6" block right
#    before the (still unrenamed) "Number of Tables:" label, reusing that
#    label run's bold formatting via FormattedText so the new label gets the
#    same <w:b/><w:bCs/> run properties as the other labels.
# ---------------------------------------------------------------------------
$tablesLabel = $d.Content
$tablesLabel.Find.Text = "Number of Tables:"
$tablesLabel.Find.Execute() | Out-Null
$labelLen = "Number of Tables:".Length
$boldTemplate = $d.Range($tablesLabel.Start, $tablesLabel.End).FormattedText

$insertPoint = $d.Range($tablesLabel.Start, $tablesLabel.Start)
$insertPoint.FormattedText = $boldTemplate
# the inserted copy now occupies [Start, Start+labelLen); grab it and retext it
$newLabel = $d.Range($tablesLabel.Start, $tablesLabel.Start + $labelLen)
$newLabel.Text = "Number of Supplementary Figures:"

# continue building right after the new label: " " then "6" then a line break
$afterLabel = $d.Range($newLabel.End, $newLabel.End)
$afterLabel.InsertBefore(" ")
$afterValuePos = $newLabel.End + 1
$afterValue = $d.Range($afterValuePos, $afterValuePos)
$afterValue.InsertBefore("6")
$afterBreakPos = $afterValuePos + 1
$afterBreak = $d.Range($afterBreakPos, $afterBreakPos)
$lineBreak = [char]11
$afterBreak.InsertBefore($lineBreak)

# ---------------------------------------------------------------------------
# 4. "Number of Tables:" -> "Number of Main Tables:"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Number of Tables:", $true, $false, $false, $false, $false, $true, 1, $false,
    "Number of Main Tables:", 2) | Out-Null

# "5 main tables" -> "5"
$d.Content.Find.Execute("5 main tables", $true, $false, $false, $false, $false, $true, 1, $false,
    "5", 2) | Out-Null

# ---------------------------------------------------------------------------
# 5. "Number of Supplementary Files:" -> "Number of Supplementary Tables:"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Number of Supplementary Files:", $true, $false, $false, $false, $false, $true, 1, $false,
    "Number of Supplementary Tables:", 2) | Out-Null

# "3 (additional figures and data)" -> "4"
$d.Content.Find.Execute("3 (additional figures and data)", $true, $false, $false, $false, $false, $true, 1, $false,
    "4", 2) | Out-Null

# ---------------------------------------------------------------------------
# 6. References count: 8 -> 15
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("References: 8", $true, $false, $false, $false, $false, $true, 1, $false,
    "References: 15", 2) | Out-Null

Write-Host "Done:" $d.Paragraphs(12).Range.Text
